$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value changes (Map 1 - Gun Damage, Explosive, and Trap Values) ---
$ws.Range("C3").Value = 360
$ws.Range("C4").Value = 56
$ws.Range("C5").Value = 75
$ws.Range("C6").Value = 270
$ws.Range("C16").Value = 180
$ws.Range("C17").Value = 100
$ws.Range("C21").Value = 250

# --- View changes: scroll back to top (remove topLeftCell="B13") and ---
# --- move the selection to C21 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C21").Select()
